$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add column E, mirroring the formatting of column D for each row ---

# Row 3: empty cell with thick-bottom style (same style as D3)
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# Row 4: header/year cell (bold, right-aligned, thick-bottom border)
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = 2020

# Row 5: data cell
$ws.Range("D5").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = 11.5

# Row 6: data cell
$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = 2.6

# Row 7: data cell, but formatted with a dedicated "0.0" number format
$ws.Range("D7").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = 2
$ws.Range("E7").NumberFormat = "0.0"

# Row 8: footer data cell (thick-bottom border)
$ws.Range("D8").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = 0.3

$excel.CutCopyMode = $false

# --- Selection moves to B15 (matches the saved sheetView) ---
$ws.Range("B15").Select()
